$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Insert a new blank row before the existing second row (TS 4.4.5.1 - Vaakyam row)
$newRow = $t.Rows.Add($t.Rows.Item(2))
$newRow.Height = 1672 / 20

$pkgPrefix = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgSuffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$cell1Body = @'
<w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>TS 4.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>4</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>3</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>2</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>Padam</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p><w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>Padam</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> No: - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>23</w:t>
            </w:r>
          </w:p><w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve">Panchaati No. - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>8</w:t>
            </w:r>
          </w:p>
'@

$cell2Body = @'
<w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Malayalam Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
                <w:lang w:val="en-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ஹே</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:highlight w:val="yellow"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>தி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">| </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>வாத</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>†</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">| </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ப்ரஹே</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>†</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>தி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ரிதி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve">ப்ர </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">-- </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ஹே</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>தி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>|</w:t>
            </w:r>
          </w:p>
'@

$cell3Body = @'
<w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ஹே</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>தி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:highlight w:val="green"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">| </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>வாத</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>†</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">| </w:t>
            </w:r>
          </w:p><w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ப்ரஹே</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>†</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>தி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ரிதி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve">ப்ர </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve">-- </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>ஹே</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>தி</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
                <w:b/>
                <w:sz w:val="36"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>…</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t>:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:cs/>
                <w:lang w:bidi="ta-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>|</w:t>
            </w:r>
          </w:p><w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cs="BRH Malayalam Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
                <w:lang w:val="en-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>(</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>visargam</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t xml:space="preserve"> inserted</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="36"/>
                <w:szCs w:val="36"/>
              </w:rPr>
              <w:t>)</w:t>
            </w:r>
          </w:p>
'@

$t.Cell(2,1).Range.InsertXML($pkgPrefix + $cell1Body + $pkgSuffix)
$t.Cell(2,2).Range.InsertXML($pkgPrefix + $cell2Body + $pkgSuffix)
$t.Cell(2,3).Range.InsertXML($pkgPrefix + $cell3Body + $pkgSuffix)

Write-Host "Row inserted. Table1 rows now:" $t.Rows.Count
